$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25, shifting existing rows 25-27 down to 26-28
$ws.Rows.Item(25).Insert()

# Fill in the new row 25 with the new task entry
$ws.Range("A25").Value = "Victory and defeat jingles"
$ws.Range("B25").Value = "Feature"

# Match style of the new cells to the rest of the data rows (copy formatting from row 26, which was row 25 before insert)
$ws.Range("A26:B26").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)

# Update the active cell/selection to match the target state
$ws.Range("A25").Select()
